# Daily update at 8 AM UTC
# Row 25 is no longer the "latest" row, so it loses the date-only
# number format and gets the regular timestamp format used by all
# the other non-final rows. A new final row (26) is appended with
# the next day's data and the date-only number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 was previously the last row (date-only format); now that a new
# row is being appended below it, it reverts to the standard format.
$ws.Cells.Item(25, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 26.
$ws.Cells.Item(26, 1).Value = 45610
$ws.Cells.Item(26, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(26, 2).Value = 63
$ws.Cells.Item(26, 3).Value = 54
$ws.Cells.Item(26, 4).Value = 59
